$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 23; this shifts the existing
# rows 23-48 down to 24-49 (matching the rest of the diff, which is
# just that shift).
$ws.Rows.Item(23).Insert()

# Populate the newly inserted row 23 with the new asparagus price record.
$ws.Range("A23").Value = 5
$ws.Range("B23").Value = "Macroferia Regional de Talca"
$ws.Range("C23").Value = "Maule"
$ws.Range("D23").Value = 44512
$ws.Range("E23").Value = 7
$ws.Range("F23").Value = 300000000
$ws.Range("G23").Value = "Espárragos"
$ws.Range("H23").Value = "Verde"
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 5000
$ws.Range("K23").Value = 800
$ws.Range("L23").Value = 800
$ws.Range("M23").Value = 800
$ws.Range("N23").Value = "$/kilo"
$ws.Range("O23").Value = "Región del Maule"
$ws.Range("P23").Value = 800
$ws.Range("Q23").Value = 1
$ws.Range("R23").Value = "Hortaliza"
